$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data

$ws.Range('D2').Value = '61.123.63'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '2.930.08'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('E4').Value = '  +0.05%  '
$cell = $ws.Range('D5')
$cell.Value = "'591.62"
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +0.95%  '
$cell = $ws.Range('D6')
$cell.Value = "'146.22"
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.47%  '
$cell = $ws.Range('D9')
$cell.Value = "'6.89"
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  +3.29%  '
$cell = $ws.Range('D10')
$cell.Value = "'0.145"
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -0.06%  '
$ws.Range('E11').Value = '  -1.13%  '
$ws.Range('E12').Value = '  +1.30%  '
$cell = $ws.Range('D13')
$cell.Value = "'33.80"
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').Value = '3.413.58'
$ws.Range('E15').Value = '  +1.11%  '
$ws.Range('D16').Value = '61.082.73'
$ws.Range('E16').Value = '  +1.13%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$cell = $ws.Range('D17')
$cell.Value = "'6.71"
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -1.39%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.926.86'
$ws.Range('E18').Value = '  +1.17%  '
$cell = $ws.Range('D19')
$cell.Value = "'431.80"
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  +1.51%  '
$ws.Range('E20').Value = '  -1.22%  '
$cell = $ws.Range('D21')
$cell.Value = "'0.684"
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +2.14%  '
$ws.Range('E22').Value = '  -0.07%  '
$cell = $ws.Range('D23')
$cell.Value = "'81.42"
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  +0.76%  '
$cell = $ws.Range('D24')
$cell.Value = "'11.05"
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  +0.98%  '
$cell = $ws.Range('D25')
$cell.Value = "'2.24"
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +2.41%  '
$cell = $ws.Range('D26')
$cell.Value = "'12.05"
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +2.46%  '
$ws.Range('E28').Value = '  +6.98%  '
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  +0.17%  '
$ws.Range('E31').Value = '  -2.64%  '
$cell = $ws.Range('D32')
$cell.Value = "'26.51"
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  +0.22%  '
$ws.Range('E33').Value = '  +2.03%  '
$ws.Range('E34').Value = '  +3.19%  '
$ws.Range('E35').Value = '  +0.80%  '
$cell = $ws.Range('D36')
$cell.Value = "'5.63"
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  -0.57%  '
$cell = $ws.Range('D37')
$cell.Value = "'3.08"
$cell.Style = 'Normal'
$ws.Range('E37').Value = '  +4.25%  '
$cell = $ws.Range('D38')
$cell.Value = "'50.01"
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  +1.56%  '
$ws.Range('E39').Value = '  +2.57%  '
$ws.Range('E40').Value = '  -1.15%  '
$cell = $ws.Range('D41')
$cell.Value = "'8.60"
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -1.60%  '
$cell = $ws.Range('D42')
$cell.Value = "'0.288"
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  +0.05%  '
$cell = $ws.Range('D43')
$cell.Value = "'39.54"
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -4.14%  '
$cell = $ws.Range('D44')
$cell.Value = "'375.69"
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +1.37%  '
$ws.Range('E45').Value = '  +0.88%  '
$ws.Range('D46').Value = '2.714.67'
$ws.Range('E46').Value = '  +2.50%  '
$cell = $ws.Range('D47')
$cell.Value = "'131.63"
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -1.02%  '
$ws.Range('E48').Value = '  -0.02%  '
$cell = $ws.Range('D49')
$cell.Value = "'24.25"
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -2.84%  '
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('E51').Value = '  -2.49%  '
